$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 4.8

# Row 3
$ws.Range("G3").Value = 3.2
$ws.Range("H3").Value = 2.42
$ws.Range("I3").Value = 2.74
$ws.Range("J3").Value = 3.65

# Row 4
$ws.Range("G4").Value = 3.7
$ws.Range("I4").Value = 2.32
$ws.Range("J4").Value = 4

# Row 5
$ws.Range("G5").Value = 7.6
$ws.Range("H5").Value = 1.48
$ws.Range("I5").Value = 1.5
$ws.Range("J5").Value = 4.9
$ws.Range("P5").Value = 2.44
$ws.Range("Q5").Value = 1.6
$ws.Range("R5").Value = 1.57
$ws.Range("T5").Value = 1.8
$ws.Range("Z5").Value = 12
$ws.Range("AC5").Value = 12
$ws.Range("AE5").Value = 1000
$ws.Range("AG5").Value = 36
$ws.Range("AH5").Value = 28
$ws.Range("AI5").Value = 980
$ws.Range("AJ5").Value = 240
$ws.Range("AK5").Value = 110
$ws.Range("AM5").Value = 120
$ws.Range("AN5").Value = 120
$ws.Range("AO5").Value = 6.8
